$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Multivalued" column (K) is being duplicated into a new "Unique" column (L).
# Copy column K (header + the three data rows) into column L so the new column
# inherits the same cell formatting/number-format/shared-string typing as K.
$ws.Range("K4:K7").Copy($ws.Range("L4"))

# Re-label the new header cell.
$ws.Range("L4").Value2 = "Unique"

# New "Unique" column is unlocked (unlike the rest of the sheet), which is what
# produces the extra pair of style records seen in the target workbook.
$ws.Range("L4:L7").Locked = $false

# Touch the bottom-right corner of the sheet so the worksheet's used range
# (and therefore its recorded dimension) extends down to row 10 together with
# the newly added column L.
$ws.Range("L10").NumberFormat = "General"

# Match the workbook's recorded selection/active cell.
$ws.Range("E7").Select() | Out-Null
